$d = $word.ActiveDocument

$oldTitle = "Організація розгалужених процесів"
$newTitle = "Організація циклічних процесів. Ітераційні цикли"

# ---------------------------------------------------------------------
# 1) First occurrence: the «title» line on the cover page.
#    Target end state (4 runs):
#      "«"  (unchanged run)
#      " "  (brand new run, no rPr at all)
#      newTitle (run with rPr: noProof + color 000000)
#      "»"  (unchanged run)
# ---------------------------------------------------------------------
$full = $d.Content
$full.Find.Execute("«" + $oldTitle + "»")
$matchStart = $full.Start
$matchEnd = $full.End

$titleStart = $matchStart + 1
$titleEnd = $matchEnd - 1

# Rewrite just the title text (keeps it as its own run, bordered by the
# still-untouched « and » runs).
$titleRange = $d.Range($titleStart, $titleEnd)
$titleRange.Text = $newTitle

$newTitleEnd = $titleStart + $newTitle.Length

# Re-split the title away from its neighbours (a pure-formatting no-op
# touch forces a clean run boundary without merging adjacent runs that
# share identical formatting).
$newTitleRange = $d.Range($titleStart, $newTitleEnd)
$newTitleRange.Font.Bold = 1
$newTitleRange.Font.Bold = 0

# Manufacture a space run with NO run properties at all (matching a
# freshly-typed run in an untouched paragraph) by typing it into the
# scratch paragraph at the end of the story, then relocating it with
# cut/paste so it never inherits neighbouring formatting.
$endPos = $d.Content.End - 1
$scratch = $d.Range($endPos, $endPos)
$scratch.InsertAfter(" ")
$scratchRange = $d.Range($endPos, $endPos + 1)
$scratchRange.Cut()

$insertionPoint = $d.Range($titleStart, $titleStart)
$insertionPoint.Paste()

# ---------------------------------------------------------------------
# 2) Second occurrence: the "Лабораторна робота 3" heading title.
#    Simple in-place text swap (single run, keeps its own formatting).
# ---------------------------------------------------------------------
$d.Content.Find.Execute($oldTitle, $false, $false, $false, $false, $false, `
    $true, 1, $false, $newTitle, 2)

# ---------------------------------------------------------------------
# 3) Split "що менше 10" into "щ" | bookmark _GoBack | "о менше 10".
#    This also relocates the _GoBack bookmark away from the end of the
#    document (bookmark names are unique, so re-adding it removes the
#    old one automatically).
# ---------------------------------------------------------------------
$phrase = "для 0 ≤ x ≤ 2 з точністю до члена ряду, що менше 10"
$prefix = "для 0 ≤ x ≤ 2 з точністю до члена ряду, щ"
$suffix = "о менше 10"

$full2 = $d.Content
$full2.Find.Execute($phrase)
$phraseStart = $full2.Start
$phraseEnd = $full2.End
$splitPos = $phraseStart + $prefix.Length

$splitPoint = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# Force the trailing half to be rebuilt as a brand-new run (drops stale
# identity picked up from the original, now-truncated run) while ending
# up with the exact same visible text.
$tail = $d.Range($splitPos, $phraseEnd)
$tail.Text = "XXXXXXXXXX"
$tail2 = $d.Range($splitPos, $splitPos + $suffix.Length)
$tail2.Text = $suffix

Write-Output "done"
